$wb = $excel.ActiveWorkbook

# --- Sheet "consumers" (sheet1): append row 40 ---
$wsConsumers = $wb.Worksheets.Item("consumers")
$wsConsumers.Range("A40").Value = 28
$wsConsumers.Range("B40").Value = "Test1"
$wsConsumers.Range("C40").Value = "Test2"
$wsConsumers.Range("D40").Value = "Amviseri"
$wsConsumers.Range("E40").Value = "rrrr"

# --- Sheet "readings" (sheet2): append row 8 ---
$wsReadings = $wb.Worksheets.Item("readings")
$wsReadings.Range("A8").Value = 7
$wsReadings.Range("B8").Value = 2
$wsReadings.Range("C8").Value = "2025-06-17T09:06:04.199Z"
$wsReadings.Range("D8").Value = 153
$wsReadings.Range("E8").Value = 174
$wsReadings.Range("F8").Value = 21
$wsReadings.Range("G8").Value = 10.5
